$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.5187273333333333
$ws.Range("H2").Value = 1.556182
$ws.Range("I2").Value = 0.01248695061656416
$ws.Range("J2").Value = 0.01248695061656416
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.777551666666666
$ws.Range("N2").Value = 26.332655
$ws.Range("O2").Value = 0.07454818073713242
$ws.Range("P2").Value = 0.07454818073713242
$ws.Range("Q2").Value = 4.553155969245555
$ws.Range("R2").Value = 40.97840372321
$ws.Range("S2").Value = 0.0009308794514192719
$ws.Range("T2").Value = 0.000930879451419272

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5187273333333333
$ws.Range("H3").Value = 1.556182
$ws.Range("I3").Value = 0.01248695061656416
$ws.Range("J3").Value = 0.01248695061656416
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 36.14140700000001
$ws.Range("N3").Value = 108.424221
$ws.Range("O3").Value = 0.306950758417288
$ws.Range("P3").Value = 0.306950758417288
$ws.Range("Q3").Value = 18.74753567602467
$ws.Range("R3").Value = 168.727821084222
$ws.Range("S3").Value = 0.003832878962073591
$ws.Range("T3").Value = 0.00383287896207359

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.5187273333333333
$ws.Range("H4").Value = 1.556182
$ws.Range("I4").Value = 0.01248695061656416
$ws.Range("J4").Value = 0.01248695061656416
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.580447
$ws.Range("N4").Value = 163.741341
$ws.Range("O4").Value = 0.4635544377507104
$ws.Range("P4").Value = 0.4635544377507104
$ws.Range("Q4").Value = 28.31236972445133
$ws.Range("R4").Value = 254.811327520062
$ws.Range("S4").Value = 0.005788381372282284
$ws.Range("T4").Value = 0.005788381372282285

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5187273333333333
$ws.Range("H5").Value = 1.556182
$ws.Range("I5").Value = 0.01248695061656416
$ws.Range("J5").Value = 0.01248695061656416
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.24393266666667
$ws.Range("N5").Value = 54.731798
$ws.Range("O5").Value = 0.1549466230948692
$ws.Range("P5").Value = 0.1549466230948692
$ws.Range("Q5").Value = 9.463626541692889
$ws.Range("R5").Value = 85.172638875236
$ws.Range("S5").Value = 0.001934810830789012
$ws.Range("T5").Value = 0.001934810830789011

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.618566
$ws.Range("H6").Value = 58.855698
$ws.Range("I6").Value = 0.4722636519567852
$ws.Range("J6").Value = 0.4722636519567853
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.777551666666666
$ws.Range("N6").Value = 26.332655
$ws.Range("O6").Value = 0.07454818073713242
$ws.Range("P6").Value = 0.07454818073713242
$ws.Range("Q6").Value = 172.20297669091
$ws.Range("R6").Value = 1549.82679021819
$ws.Range("S6").Value = 0.03520639608165262
$ws.Range("T6").Value = 0.03520639608165263

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.618566
$ws.Range("H7").Value = 58.855698
$ws.Range("I7").Value = 0.4722636519567852
$ws.Range("J7").Value = 0.4722636519567853
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 36.14140700000001
$ws.Range("N7").Value = 108.424221
$ws.Range("O7").Value = 0.306950758417288
$ws.Range("P7").Value = 0.306950758417288
$ws.Range("Q7").Value = 709.0425785623622
$ws.Range("R7").Value = 6381.38320706126
$ws.Range("S7").Value = 0.1449616861410534
$ws.Range("T7").Value = 0.1449616861410534

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.618566
$ws.Range("H8").Value = 58.855698
$ws.Range("I8").Value = 0.4722636519567852
$ws.Range("J8").Value = 0.4722636519567853
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.580447
$ws.Range("N8").Value = 163.741341
$ws.Range("O8").Value = 0.4635544377507104
$ws.Range("P8").Value = 0.4635544377507104
$ws.Range("Q8").Value = 1070.790101779002
$ws.Range("R8").Value = 9637.110916011019
$ws.Range("S8").Value = 0.2189199116529247
$ws.Range("T8").Value = 0.2189199116529248

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.618566
$ws.Range("H9").Value = 58.855698
$ws.Range("I9").Value = 0.4722636519567852
$ws.Range("J9").Value = 0.4722636519567853
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.24393266666667
$ws.Range("N9").Value = 54.731798
$ws.Range("O9").Value = 0.1549466230948692
$ws.Range("P9").Value = 0.1549466230948692
$ws.Range("Q9").Value = 357.9197971205561
$ws.Range("R9").Value = 3221.278174085005
$ws.Range("S9").Value = 0.0731756580811545
$ws.Range("T9").Value = 0.07317565808115449

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.284536
$ws.Range("H10").Value = 57.85360799999999
$ws.Range("I10").Value = 0.4642227876212815
$ws.Range("J10").Value = 0.4642227876212815
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.777551666666666
$ws.Range("N10").Value = 26.332655
$ws.Range("O10").Value = 0.07454818073713242
$ws.Range("P10").Value = 0.07454818073713242
$ws.Range("Q10").Value = 169.2710111076933
$ws.Range("R10").Value = 1523.43909996924
$ws.Range("S10").Value = 0.03460696427388673
$ws.Range("T10").Value = 0.03460696427388673

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 19.284536
$ws.Range("H11").Value = 57.85360799999999
$ws.Range("I11").Value = 0.4642227876212815
$ws.Range("J11").Value = 0.4642227876212815
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 36.14140700000001
$ws.Range("N11").Value = 108.424221
$ws.Range("O11").Value = 0.306950758417288
$ws.Range("P11").Value = 0.306950758417288
$ws.Range("Q11").Value = 696.9702643821521
$ws.Range("R11").Value = 6272.732379439369
$ws.Range("S11").Value = 0.14249353673494
$ws.Range("T11").Value = 0.1424935367349399

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 19.284536
$ws.Range("H12").Value = 57.85360799999999
$ws.Range("I12").Value = 0.4642227876212815
$ws.Range("J12").Value = 0.4642227876212815
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.580447
$ws.Range("N12").Value = 163.741341
$ws.Range("O12").Value = 0.4635544377507104
$ws.Range("P12").Value = 0.4635544377507104
$ws.Range("Q12").Value = 1052.558595067592
$ws.Range("R12").Value = 9473.027355608327
$ws.Range("S12").Value = 0.2151925333068506
$ws.Range("T12").Value = 0.2151925333068506

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 19.284536
$ws.Range("H13").Value = 57.85360799999999
$ws.Range("I13").Value = 0.4642227876212815
$ws.Range("J13").Value = 0.4642227876212815
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.24393266666667
$ws.Range("N13").Value = 54.731798
$ws.Range("O13").Value = 0.1549466230948692
$ws.Range("P13").Value = 0.1549466230948692
$ws.Range("Q13").Value = 351.8257762919094
$ws.Range("R13").Value = 3166.431986627184
$ws.Range("S13").Value = 0.07192975330560424
$ws.Range("T13").Value = 0.07192975330560422

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.119724666666667
$ws.Range("H14").Value = 6.359174
$ws.Range("I14").Value = 0.05102660980536902
$ws.Range("J14").Value = 0.05102660980536902
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.777551666666666
$ws.Range("N14").Value = 26.332655
$ws.Range("O14").Value = 0.07454818073713242
$ws.Range("P14").Value = 0.07454818073713242
$ws.Range("Q14").Value = 18.60599278077445
$ws.Range("R14").Value = 167.45393502697
$ws.Range("S14").Value = 0.003803940930173783
$ws.Range("T14").Value = 0.003803940930173783

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.119724666666667
$ws.Range("H15").Value = 6.359174
$ws.Range("I15").Value = 0.05102660980536902
$ws.Range("J15").Value = 0.05102660980536902
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 36.14140700000001
$ws.Range("N15").Value = 108.424221
$ws.Range("O15").Value = 0.306950758417288
$ws.Range("P15").Value = 0.306950758417288
$ws.Range("Q15").Value = 76.60983190593936
$ws.Range("R15").Value = 689.4884871534541
$ws.Range("S15").Value = 0.01566265657922105
$ws.Range("T15").Value = 0.01566265657922104

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.119724666666667
$ws.Range("H16").Value = 6.359174
$ws.Range("I16").Value = 0.05102660980536902
$ws.Range("J16").Value = 0.05102660980536902
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.580447
$ws.Range("N16").Value = 163.741341
$ws.Range("O16").Value = 0.4635544377507104
$ws.Range("P16").Value = 0.4635544377507104
$ws.Range("Q16").Value = 115.6955198235927
$ws.Range("R16").Value = 1041.259678412334
$ws.Range("S16").Value = 0.02365361141865272
$ws.Range("T16").Value = 0.02365361141865272

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.119724666666667
$ws.Range("H17").Value = 6.359174
$ws.Range("I17").Value = 0.05102660980536902
$ws.Range("J17").Value = 0.05102660980536902
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.24393266666667
$ws.Range("N17").Value = 54.731798
$ws.Range("O17").Value = 0.1549466230948692
$ws.Range("P17").Value = 0.1549466230948692
$ws.Range("Q17").Value = 38.67211409053913
$ws.Range("R17").Value = 348.049026814852
$ws.Range("S17").Value = 0.007906400877321472
$ws.Range("T17").Value = 0.007906400877321472
